$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect it so the cell values/text can be updated.
$ws.Unprotect()

# Update the confidential disclaimer text (date change 2021-05-12 -> 2021-05-13)
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-13 for illustrative purposes only and are subject to change."

# Update the numeric values in D2:E4 (weights / percent changes)
$ws.Range("D2").Value = 0.8483147294611066
$ws.Range("E2").Value = 0.007467490665636678

$ws.Range("D3").Value = 0.1516852705388934
$ws.Range("E3").Value = -0.001344086021505375

$ws.Range("E4").Value = 0.00613090427197327

# Restore sheet protection (contents protected) as it was before the edit.
$ws.Protect()
